# Consolidating 24Q1 through 25Q3 data into 25Q3 report.
# Drop the oldest quarter (Q2 2024) column, shift the remaining quarter
# columns left by one, and append the newest quarter (Q2 2025, N=46) data
# in the last column. Three rows (Phone / Door hangers / Facebook...) are
# also reordered to match the source's new ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: quarter labels shift left by one, newest quarter added ---
$ws.Range("B1").Value = "Q3 2024" + [char]10 + "(N=19)"
$ws.Range("C1").Value = "Q4 2024" + [char]10 + "(N=33)"
$ws.Range("D1").Value = "Q1 2025" + [char]10 + "(N=45)"
$ws.Range("E1").Value = "Q2 2025" + [char]10 + "(N=46)"
# Recompute the row height so it isn't pinned to an explicit "custom"
# height after the multi-line header text is (re)written.
$ws.Rows("1").AutoFit()

# --- Row 2: Email ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.7111111111111111
$ws.Range("E2").Value = 0.5869565217391305

# --- Row 3: Direct mail ---
$ws.Range("B3").Value = 0.1052631578947368
$ws.Range("C3").Value = 0.09090909090909093
$ws.Range("D3").Value = 0.4444444444444444
$ws.Range("E3").Value = 0.4130434782608696

# --- Row 4: Text messages ---
$ws.Range("B4").Value = 0.2105263157894737
$ws.Range("C4").Value = 0.4545454545454545
$ws.Range("D4").Value = 0.3777777777777778
$ws.Range("E4").Value = 0.3043478260869565

# --- Row 5: Austin Energy's website ---
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.3555555555555556
$ws.Range("E5").Value = 0.2173913043478261

# --- Row 6: was "Door hangers" -> now "Phone" ---
$ws.Range("A6").Value = "Phone"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0.06666666666666667
$ws.Range("E6").Value = 0.1521739130434783

# --- Row 7: was "Facebook, Twitter, or other Social Media" -> now "Door hangers" ---
$ws.Range("A7").Value = "Door hangers"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.1333333333333333
$ws.Range("E7").Value = 0.06521739130434782

# --- Row 8: was "Phone" -> now "Facebook, Twitter, or other Social Media" ---
$ws.Range("A8").Value = "Facebook, Twitter, or other Social Media"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0.0303030303030303
$ws.Range("D8").Value = 0.08888888888888889
$ws.Range("E8").Value = 0.06521739130434782

# --- Row 9: Email (trailing space in label, unchanged) ---
$ws.Range("B9").Value = 0.631578947368421
$ws.Range("C9").Value = 0.7575757575757576
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0

# --- Row 10: Phone call ---
$ws.Range("B10").Value = 0.3157894736842105
$ws.Range("C10").Value = 0.2424242424242424
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0

# --- Row 11: Austin Energy's Website ---
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0.1212121212121212
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0

# --- Row 12: All other ---
$ws.Range("B12").Value = 0.05263157894736842
$ws.Range("C12").Value = 0.05263157894736842
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.06521739130434782
